$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 22 (AssignQualification) gets a DONE status in column E (Column1)
$ws.Range("E22").Value = "DONE"

# Row 23 (GrantQualification) gets IN PROGRESS status in column E (Column1)
# and a note in column F (Column2)
$ws.Range("E23").Value = "IN PROGRESS"
$ws.Range("F23").Value = "REQUIRES Qualification Request entity creation"

# Move the active selection to F23 to match the author's last edited cell
$ws.Range("F23").Select()
